$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.148840308189392
$ws.Range("B1").Value = 1.332544922828674
$ws.Range("C1").Value = 1.691295981407166
$ws.Range("D1").Value = 3.213660478591919
$ws.Range("E1").Value = -1
